$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, pushing the existing "Wonthaggi" row
# (and everything after it) down to row 54.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the "Southern Cross" entry.
$ws.Range("A53").Value = "Southern Cross"
$ws.Range("B53").Value = "Metro trains - Mernda line"
$ws.Range("C53").Value = "28/12/2020 14:30 - 14:45"
$ws.Range("D53").Value = "Caught train from Southern Cross to Victoria Park station"
